$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataEntry")

# Update the label in A5 to reflect the new wording.
$ws.Range("A5").Value = "SensorEffector choice"

# The cell wraps text, so the longer wording now spans two lines;
# adjust the row height so it matches the new, taller content.
$ws.Rows.Item(5).RowHeight = 30

# Move the active selection to A5 to match the saved view state.
$ws.Range("A5").Select()
